$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new "category" column before the existing "date" column (I),
# pushing date/legislator_name/legislator_id one column to the right.
$ws.Columns("I:I").Insert()
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(2, 9).Value = "normal"

# Append two new trailing columns: source_file and index.
# Copy formatting from the neighbouring header/data cells first so the new
# cells pick up the same style (bold/centered header vs plain data row).
$ws.Cells.Item(1, 11).Copy($ws.Cells.Item(1, 13))
$ws.Cells.Item(1, 11).Copy($ws.Cells.Item(1, 14))
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(2, 13))
$ws.Cells.Item(2, 11).Copy($ws.Cells.Item(2, 14))

$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"
$ws.Cells.Item(2, 13).Value = "tmp399c1"
$ws.Cells.Item(2, 14).Value = 63
